# Insert a new "date" column before the existing "in" column (which was
# column D) so the attendance template gains a date field between
# employee_code and in/out/status.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift in/out/status one column to the right, opening up D1 for "date".
$ws.Columns.Item(4).Insert()

$ws.Range("D1").Value = "date"

# Mirror the "Ignore Error" action a user takes in Excel to keep the
# number-stored-as-text warning suppressed across the (now wider) header row.
$ws.Range("A1:G1").Errors.Item(9).Ignore = $true
